# Update the salinity calibration results with the new training run values
# (training_test_2016_2019, k=10) while keeping the same row/column headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged labels)
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "RMSE"
$ws.Range("C1").Value = "NRMSE"
$ws.Range("D1").Value = "MAE"
$ws.Range("E1").Value = "RSE"
$ws.Range("F1").Value = "RRSE"
$ws.Range("G1").Value = "RAE"
$ws.Range("H1").Value = "R2"
$ws.Range("I1").Value = "Corr Coeff"

# Row labels (unchanged)
$ws.Range("A2").Value = "random_forest"
$ws.Range("A3").Value = "lsboost"
$ws.Range("A4").Value = "neural_network"

# random_forest row - updated metrics
$ws.Range("B2").Value = 3.6123774410934022
$ws.Range("C2").Value = 0.25803546167868258
$ws.Range("D2").Value = 2.6685154490573511
$ws.Range("E2").Value = 0.25889930797023619
$ws.Range("F2").Value = 0.50882148929682225
$ws.Range("G2").Value = 0.4539140156098056
$ws.Range("H2").Value = 0.74110069202976381
$ws.Range("I2").Value = 0.86088397525680527

# lsboost row - updated metrics
$ws.Range("B3").Value = 3.6765139028616494
$ws.Range("C3").Value = 0.26261678846213121
$ws.Range("D3").Value = 2.7507392917803912
$ws.Range("E3").Value = 0.26817424821395264
$ws.Range("F3").Value = 0.51785543177025062
$ws.Range("G3").Value = 0.46790027701311443
$ws.Range("H3").Value = 0.73182575178604736
$ws.Range("I3").Value = 0.85554979588220958

# neural_network row - updated metrics
$ws.Range("B4").Value = 3.7571622577195241
$ws.Range("C4").Value = 0.2683775750407007
$ws.Range("D4").Value = 2.7726271779127898
$ws.Range("E4").Value = 0.28006868265378426
$ws.Range("F4").Value = 0.52921515724115864
$ws.Range("G4").Value = 0.471623402652532
$ws.Range("H4").Value = 0.71993131734621574
$ws.Range("I4").Value = 0.84868114445409182
